$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17 (anchor G17=38956)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 372081.06
$ws.Range("J17").Value = 372081.06
$ws.Range("L17").Value = 1116243.18
$ws.Range("N17").Value = -1116579.18

# Sheet ALC, row 38 (anchor G38=4599)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1440.5
$ws.Range("I38").Value = 45.11111
$ws.Range("K38").Value = 135.33333
$ws.Range("M38").Value = 236.66667

# Sheet ALC, row 61 (anchor G61=4604)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 4823
$ws.Range("I61").Value = 4823
$ws.Range("K61").Value = 14469
$ws.Range("M61").Value = -14297

# Sheet ALC, row 64 (anchor G64=5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3970.3333
$ws.Range("J64").Value = 6330
$ws.Range("L64").Value = 6330
$ws.Range("N64").Value = -6826

# Sheet ALC, row 67 (anchor G67=5506)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3970.3333
$ws.Range("J67").Value = 6330
$ws.Range("L67").Value = 6330
$ws.Range("N67").Value = -8046

# Sheet ALC, row 88 (anchor G88=12608)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1118.5
$ws.Range("I88").Value = 1051.8
$ws.Range("J88").Value = 1201.875
$ws.Range("K88").Value = 1051.8
$ws.Range("L88").Value = 1201.875
$ws.Range("M88").Value = -645.8
$ws.Range("N88").Value = -2013.875

# Sheet ALC, row 91 (anchor G91=12608)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 1118.5
$ws.Range("I91").Value = 1051.8
$ws.Range("J91").Value = 1201.875
$ws.Range("K91").Value = 1051.8
$ws.Range("L91").Value = 1201.875
$ws.Range("M91").Value = 352.2
$ws.Range("N91").Value = -4009.875

# Sheet ALC, row 96 (anchor G96=19894)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2324
$ws.Range("I96").Value = 2417.75
$ws.Range("J96").Value = 2199
$ws.Range("K96").Value = 7253.25
$ws.Range("L96").Value = 6597
$ws.Range("M96").Value = -5880.25
$ws.Range("N96").Value = -9343

# Sheet ALC, row 132 (anchor G132=44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1857.8368
$ws.Range("I132").Value = 1296.3636
$ws.Range("K132").Value = 3889.0908
$ws.Range("M132").Value = -1359.0908

# Sheet ALC, row 137 (anchor G137=44013)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 29189.854
$ws.Range("I137").Value = 38485.9
$ws.Range("J137").Value = 3837
$ws.Range("K137").Value = 115457.7
$ws.Range("L137").Value = 11511
$ws.Range("M137").Value = -112907.7
$ws.Range("N137").Value = -16611

# Sheet ARM, row 45 (anchor G45=27714)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1821.8889
$ws.Range("I45").Value = 1232.8334
$ws.Range("K45").Value = 1232.8334
$ws.Range("M45").Value = -855.8334

# Sheet ARM, row 122 (anchor G122=36168)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2819.647
$ws.Range("I122").Value = 2883.625
$ws.Range("K122").Value = 8650.875
$ws.Range("M122").Value = -6200.875

# Sheet ARM, row 132 (anchor G132=43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1608.8113
$ws.Range("I132").Value = 1537.1666
$ws.Range("J132").Value = 2296.6
$ws.Range("K132").Value = 4611.4998
$ws.Range("L132").Value = 6889.799999999999
$ws.Range("M132").Value = -2081.4998
$ws.Range("N132").Value = -11949.8

# Sheet BSM, row 94 (anchor G94=19939)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1842.5385
$ws.Range("I94").Value = 1789.4445
$ws.Range("K94").Value = 1789.4445
$ws.Range("M94").Value = -1338.4445

# Sheet BSM, row 105 (anchor G105=19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2415.8572
$ws.Range("I105").Value = 2501.8333
$ws.Range("K105").Value = 2501.8333
$ws.Range("M105").Value = -754.8332999999998

# Sheet BSM, row 134 (anchor G134=43998)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2369.075
$ws.Range("I134").Value = 1996.2572
$ws.Range("K134").Value = 5988.7716
$ws.Range("M134").Value = -3453.7716

# Sheet CRP, row 31 (anchor G31=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 324021.6
$ws.Range("I31").Value = 371550.72
$ws.Range("J31").Value = 3200
$ws.Range("K31").Value = 371550.72
$ws.Range("L31").Value = 3200
$ws.Range("M31").Value = -371255.72
$ws.Range("N31").Value = -3790

# Sheet CRP, row 34 (anchor G34=44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 324021.6
$ws.Range("I34").Value = 371550.72
$ws.Range("J34").Value = 3200
$ws.Range("K34").Value = 371550.72
$ws.Range("L34").Value = 3200
$ws.Range("M34").Value = -371348.72
$ws.Range("N34").Value = -3604

# Sheet CUL, row 37 (anchor G37=9516)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 166758320
$ws.Range("J37").Value = 166758320
$ws.Range("L37").Value = 500274960
$ws.Range("N37").Value = -500275184

# Sheet GSM, row 102 (anchor G102=36169)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3119.3845
$ws.Range("J102").Value = 2301
$ws.Range("L102").Value = 2301
$ws.Range("N102").Value = -5545

# Sheet GSM, row 122 (anchor G122=36182)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2943367.5
$ws.Range("I122").Value = 4167899.8
$ws.Range("J122").Value = 4490
$ws.Range("K122").Value = 12503699.4
$ws.Range("L122").Value = 13470
$ws.Range("M122").Value = -12501249.4
$ws.Range("N122").Value = -18370

# Sheet GSM, row 132 (anchor G132=44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 30832.578
$ws.Range("I132").Value = 32400.084
$ws.Range("J132").Value = 2617.5
$ws.Range("K132").Value = 97200.25199999999
$ws.Range("L132").Value = 7852.5
$ws.Range("M132").Value = -94670.25199999999
$ws.Range("N132").Value = -12912.5

# Sheet GSM, row 136 (anchor G136=42218)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 21419
$ws.Range("J136").Value = 21419
$ws.Range("L136").Value = 64257
$ws.Range("N136").Value = -69357

# Sheet LTW, row 20 (anchor G20=4308)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 10654.667
$ws.Range("I20").Value = 2000
$ws.Range("J20").Value = 14982
$ws.Range("K20").Value = 2000
$ws.Range("L20").Value = 14982
$ws.Range("M20").Value = -1774
$ws.Range("N20").Value = -15434

# Sheet LTW, row 46 (anchor G46=5282)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2745.5454
$ws.Range("I46").Value = 1689.1666
$ws.Range("K46").Value = 1689.1666
$ws.Range("M46").Value = -1501.1666

# Sheet LTW, row 136 (anchor G136=44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1721.8049
$ws.Range("I136").Value = 1167.303
$ws.Range("J136").Value = 4009.125
$ws.Range("K136").Value = 3501.909000000001
$ws.Range("L136").Value = 12027.375
$ws.Range("M136").Value = -951.9090000000006
$ws.Range("N136").Value = -17127.375

# Sheet WVR, row 107 (anchor G107=27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1663.1111
$ws.Range("I107").Value = 1599.7142
$ws.Range("J107").Value = 1885
$ws.Range("K107").Value = 4799.142599999999
$ws.Range("L107").Value = 5655
$ws.Range("M107").Value = -2879.142599999999
$ws.Range("N107").Value = -9495

# Sheet WVR, row 122 (anchor G122=36208)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2840.3076
$ws.Range("I122").Value = 2328.889
$ws.Range("J122").Value = 3991
$ws.Range("K122").Value = 6986.667
$ws.Range("L122").Value = 11973
$ws.Range("M122").Value = -4536.667
$ws.Range("N122").Value = -16873

# Sheet WVR, row 132 (anchor G132=44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1042
$ws.Range("J132").Value = 1149.5
$ws.Range("L132").Value = 3448.5
$ws.Range("N132").Value = -8508.5

# Sheet WVR, row 136 (anchor G136=44031)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 211743.38
$ws.Range("I136").Value = 253374.67
$ws.Range("J136").Value = 3586.875
$ws.Range("K136").Value = 760124.01
$ws.Range("L136").Value = 10760.625
$ws.Range("M136").Value = -757574.01
$ws.Range("N136").Value = -15860.625
